# Updated cryptos list (price + 1h volume refresh) as produced by the
# scheduled "Updated cryptos list ... with GitHub Actions" commit.
# A few coins also swapped rank/row position (Monero/EthereumClassic,
# Filecoin/SuiNetwork, FirstDigitalUSD/EnergySwap) so their B (name) and
# C (link) columns are rewritten too, not just D (price) / E (volume).
#
# Price column D is stored as plain text in the workbook (e.g. "6.50",
# "0.840") even though it looks numeric, so that trailing zeros survive.
# For the handful of new prices that are themselves "clean" decimals
# (6.50, 0.840, 0.990) we force the cell to Text format first so Excel
# doesn't silently renumber them to 6.5 / 0.84 / 0.99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.584.43"
$ws.Range("E2").Value = "  +6.54%  "

$ws.Range("D3").Value = "2.645.64"
$ws.Range("E3").Value = "  +8.80%  "

$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.65%  "

$ws.Range("D5").Value = "512.32"
$ws.Range("E5").Value = "  +4.86%  "

$ws.Range("D6").Value = "157.83"
$ws.Range("E6").Value = "  +2.55%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").Value = "2.679.71"
$ws.Range("E9").Value = "  +10.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.50"
$ws.Range("E10").Value = "  +2.72%  "

$ws.Range("E11").Value = "  +5.07%  "

$ws.Range("D12").Value = "0.349"
$ws.Range("E12").Value = "  +3.91%  "

$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").Value = "3.149.24"
$ws.Range("E14").Value = "  +11.19%  "

$ws.Range("D15").Value = "60.760.53"
$ws.Range("E15").Value = "  +6.76%  "

$ws.Range("D16").Value = "21.85"
$ws.Range("E16").Value = "  +5.44%  "

$ws.Range("E17").Value = "  +5.18%  "

$ws.Range("D18").Value = "2.686.53"
$ws.Range("E18").Value = "  +10.75%  "

$ws.Range("D19").Value = "4.81"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("D20").Value = "349.67"
$ws.Range("E20").Value = "  +7.94%  "

$ws.Range("D21").Value = "10.55"
$ws.Range("E21").Value = "  +5.34%  "

$ws.Range("D22").Value = "6.21"
$ws.Range("E22").Value = "  +3.56%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Value = "60.49"
$ws.Range("E24").Value = "  +3.98%  "

$ws.Range("E25").Value = "  +3.97%  "

$ws.Range("D26").Value = "2.801.79"
$ws.Range("E26").Value = "  +11.12%  "

$ws.Range("E27").Value = "  +3.02%  "

$ws.Range("E28").Value = "  +0.42%  "

$ws.Range("D29").Value = "0.0₃0867"
$ws.Range("E29").Value = "  +10.62%  "

$ws.Range("D30").Value = "7.54"
$ws.Range("E30").Value = "  +3.02%  "

$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "19.69"
$ws.Range("E32").Value = "  +5.98%  "

$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "157.47"
$ws.Range("E33").Value = "  +4.69%  "

$ws.Range("E34").Value = "  +3.59%  "

$ws.Range("D35").Value = "5.74"
$ws.Range("E35").Value = "  +8.60%  "

$ws.Range("D36").Value = "4.09"
$ws.Range("E36").Value = "  +9.60%  "

$ws.Range("E37").Value = "  +5.24%  "

$ws.Range("D38").Value = "314.58"
$ws.Range("E38").Value = "  +17.69%  "

$ws.Range("E39").Value = "  +10.44%  "

$ws.Range("D40").Value = "0.857"
$ws.Range("E40").Value = "  +1.19%  "

$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.840"
$ws.Range("E41").Value = "  +31.05%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "3.78"
$ws.Range("E42").Value = "  +6.75%  "

$ws.Range("D43").Value = "35.43"
$ws.Range("E43").Value = "  +3.66%  "

$ws.Range("D44").Value = "0.644"
$ws.Range("E44").Value = "  +8.55%  "

$ws.Range("D45").Value = "0.0577"
$ws.Range("E45").Value = "  +8.27%  "

$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "19.94"
$ws.Range("E47").Value = "  +14.30%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.990"
$ws.Range("E48").Value = "  -0.70%  "

$ws.Range("D49").Value = "2.073.65"
$ws.Range("E49").Value = "  +10.66%  "

$ws.Range("E50").Value = "  +3.57%  "

$ws.Range("D51").Value = "4.86"
$ws.Range("E51").Value = "  +4.87%  "
